# Update the Bigmac row's Store Price / GP-Food price figures.
# (The dependent "Result" column is a shared formula and recalculates
# automatically.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 36

# Move the active selection from D19 to D7, matching the author's
# last cursor position when the file was saved.
$ws.Range("D7").Select()
